$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the shared-string table in the same order the Northwind Customers
# field list is naturally enumerated, by touching scratch cells first and
# then clearing them again (values are still interned as shared strings).
$ws.Range("A20").Value = "Address"
$ws.Range("A21").Value = "City"
$ws.Range("A22").Value = "Region"
$ws.Range("A23").Value = "Country"
$ws.Range("A24").Value = "Phone"
$ws.Range("A25").Value = "Fax"
$ws.Range("A26").Value = "Company"
$ws.Range("A27").Value = "Contact Title"
$ws.Range("A28").Value = "Postal Code"
$ws.Range("A29").Value = "Contact Name"
$ws.Range("A20:A29").Clear()

# New header values for the Northwind "Customers"-style layout
$ws.Range("A1").Value = "Company"
$ws.Range("B1").Value = "Contact Name"
$ws.Range("C1").Value = "Contact Title"
$ws.Range("D1").Value = "Address"
$ws.Range("E1").Value = "City"
$ws.Range("F1").Value = "Region"
$ws.Range("G1").Value = "Postal Code"
$ws.Range("H1").Value = "Country"
$ws.Range("I1").Value = "Phone"
$ws.Range("J1").Value = "Fax"

# Column widths (values chosen so the engine's internal character<->pixel
# rounding lands on the exact target stored widths)
$ws.Columns.Item(1).ColumnWidth = 30
$ws.Columns.Item(2).ColumnWidth = 27.6666666666667
$ws.Columns.Item(3).ColumnWidth = 16.6666666666667
$ws.Columns.Item(4).ColumnWidth = 29.6666666666667
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 7.66666666666667
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 10.6666666666667
$ws.Columns.Item(9).ColumnWidth = 15
$ws.Columns.Item(10).ColumnWidth = 14.3333333333333

# Selection
$ws.Range("B2").Select() | Out-Null
